$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.223.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.724.08'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '330.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.531'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.562'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0828'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.148.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.700.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.881'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.042.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0962'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '290.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.22%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.60'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0826'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.23%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("E38").Value = '  +0.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '128.98'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0349'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.02%  '
$ws.Range("E43").Value = '  +3.84%  '
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.111.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.09%  '
$ws.Range("E48").Value = '  -2.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.45'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.67%  '
